# Oklahoma overview workbook - convert numeric "count" cells to text cells
# (matches the upstream factsheet pipeline change: COMM text-formatted exports)
# and append the missing "Total" row to the County sheet.

function Set-TextValue {
    param(
        $Range,
        [string]$Value
    )
    # Force a Text number format before assigning so the engine stores the
    # cell as a string (t="inlineStr"/"s") instead of re-parsing "906" etc.
    # back into a number, then drop back to the default "Normal" style so we
    # don't leave a stray numFmt applied to the cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 906 (number) -> "906" (text)
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "906"

# ---------------------------------------------------------------------
# Sheet "County": column B (rows 2-67) number -> text (same values);
# rows 68-74 (zero counties) get new text values across B:F;
# append new row 75 "Total".
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2 = "2";   3 = "2";   4 = "1";   5 = "5";   6 = "1";   7 = "10";  8 = "6";
    9 = "7";   10 = "23"; 11 = "15"; 12 = "5";  13 = "43"; 14 = "2";  15 = "22";
    16 = "6";  17 = "6";  18 = "7";  19 = "6";  20 = "1";  21 = "16"; 22 = "3";
    23 = "8";  24 = "1";  25 = "1";  26 = "2";  27 = "3";  28 = "6";  29 = "3";
    30 = "12"; 31 = "2";  32 = "2";  33 = "4";  34 = "1";  35 = "7";  36 = "2";
    37 = "2";  38 = "7";  39 = "3";  40 = "6";  41 = "2";  42 = "1";  43 = "19";
    44 = "1";  45 = "3";  46 = "1";  47 = "267"; 48 = "4"; 49 = "4";  50 = "6";
    51 = "3";  52 = "20"; 53 = "10"; 54 = "18"; 55 = "17"; 56 = "1"; 57 = "12";
    58 = "7";  59 = "3";  60 = "9";  61 = "5";  62 = "2";  63 = "201"; 64 = "2";
    65 = "18"; 66 = "3";  67 = "6";
}

foreach ($row in $countyCounts.Keys) {
    Set-TextValue $wsCounty.Cells.Item($row, 2) $countyCounts[$row]
}

# Rows 68-74: the seven counties with no filers, now shown with % / $ formatting
foreach ($row in 68..74) {
    Set-TextValue $wsCounty.Cells.Item($row, 2) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($row, 3) "`$0"
    Set-TextValue $wsCounty.Cells.Item($row, 4) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($row, 5) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($row, 6) "0.00%"
}

# New "Total" row (row 75)
Set-TextValue $wsCounty.Cells.Item(75, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(75, 2) "906"
Set-TextValue $wsCounty.Cells.Item(75, 3) "`$1,608,045,127"
Set-TextValue $wsCounty.Cells.Item(75, 4) "6.01%"
Set-TextValue $wsCounty.Cells.Item(75, 5) "-20.07%"
Set-TextValue $wsCounty.Cells.Item(75, 6) "72.63%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": column B (rows 2-7) number -> text
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")
$cdCounts = @{ 2 = "209"; 3 = "153"; 4 = "148"; 5 = "140"; 6 = "256"; 7 = "906" }
foreach ($row in $cdCounts.Keys) {
    Set-TextValue $wsCD.Cells.Item($row, 2) $cdCounts[$row]
}

# ---------------------------------------------------------------------
# Sheet "Size": column B (rows 2-8) number -> text
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @{ 2 = "260"; 3 = "274"; 4 = "166"; 5 = "68"; 6 = "92"; 7 = "46"; 8 = "906" }
foreach ($row in $sizeCounts.Keys) {
    Set-TextValue $wsSize.Cells.Item($row, 2) $sizeCounts[$row]
}

# ---------------------------------------------------------------------
# Sheet "Subsector": column B (rows 2-13) number -> text
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
$subCounts = @{
    2 = "69"; 3 = "82"; 4 = "23"; 5 = "70"; 6 = "8"; 7 = "314"; 8 = "10";
    9 = "59"; 10 = "25"; 11 = "234"; 12 = "12"; 13 = "906"
}
foreach ($row in $subCounts.Keys) {
    Set-TextValue $wsSub.Cells.Item($row, 2) $subCounts[$row]
}
